# Auto-generated edit script applying cached-value updates to the
# Ridill_Profits workbook, per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 112
$ws.Range("H112").Value = 1157.2916
$ws.Range("I112").Value = 645
$ws.Range("J112").Value = 1464.6666
$ws.Range("K112").Value = 1935
$ws.Range("L112").Value = 4393.9998
$ws.Range("M112").Value = -827
$ws.Range("N112").Value = -6609.9998
# row 137
$ws.Range("H137").Value = 17869588
$ws.Range("I137").Value = 3677241
$ws.Range("J137").Value = 78187064
$ws.Range("K137").Value = 11031723
$ws.Range("L137").Value = 234561192
$ws.Range("M137").Value = -11029173
$ws.Range("N137").Value = -234566292
# row 138
$ws.Range("H138").Value = 3125.1
$ws.Range("I138").Value = 3216.9048
$ws.Range("J138").Value = 3058.6206
$ws.Range("K138").Value = 9650.714399999999
$ws.Range("L138").Value = 9175.861800000001
$ws.Range("M138").Value = -4510.714399999999
$ws.Range("N138").Value = -19455.8618

$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 385664.66
$ws.Range("I45").Value = 667504.5600000001
$ws.Range("J45").Value = 1337.5454
$ws.Range("K45").Value = 667504.5600000001
$ws.Range("L45").Value = 1337.5454
$ws.Range("M45").Value = -667127.5600000001
$ws.Range("N45").Value = -2091.5454
# row 61
$ws.Range("H61").Value = 4568960
$ws.Range("I61").Value = 2779178.8
$ws.Range("J61").Value = 8404206
$ws.Range("K61").Value = 2779178.8
$ws.Range("L61").Value = 8404206
$ws.Range("M61").Value = -2778966.8
$ws.Range("N61").Value = -8404630
# row 74
$ws.Range("H74").Value = 128969704
$ws.Range("I74").Value = 126786310
$ws.Range("K74").Value = 126786310
$ws.Range("M74").Value = -126785436
# row 77
$ws.Range("H77").Value = 128969704
$ws.Range("I77").Value = 126786310
$ws.Range("K77").Value = 633931550
$ws.Range("M77").Value = -633927182
# row 132
$ws.Range("H132").Value = 12157083
$ws.Range("I132").Value = 12350263
$ws.Range("J132").Value = 11113911
$ws.Range("K132").Value = 37050789
$ws.Range("L132").Value = 33341733
$ws.Range("M132").Value = -37048259
$ws.Range("N132").Value = -33346793
# row 136
$ws.Range("H136").Value = 4568960
$ws.Range("I136").Value = 2779178.8
$ws.Range("J136").Value = 8404206
$ws.Range("K136").Value = 8337536.399999999
$ws.Range("L136").Value = 25212618
$ws.Range("M136").Value = -8334986.399999999
$ws.Range("N136").Value = -25217718

$ws = $wb.Worksheets.Item("BSM")
# row 105
$ws.Range("H105").Value = 1880
$ws.Range("I105").Value = 1880
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1880
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -133
$ws.Range("N105").ClearContents()
# row 134
$ws.Range("H134").Value = 24352734
$ws.Range("I134").Value = 50001680
$ws.Range("J134").Value = 2978611.5
$ws.Range("K134").Value = 150005040
$ws.Range("L134").Value = 8935834.5
$ws.Range("M134").Value = -150002505
$ws.Range("N134").Value = -8940904.5

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2483428.5
$ws.Range("I31").Value = 1390429.6
$ws.Range("J31").Value = 5215926
$ws.Range("K31").Value = 1390429.6
$ws.Range("L31").Value = 5215926
$ws.Range("M31").Value = -1390134.6
$ws.Range("N31").Value = -5216516
# row 34
$ws.Range("H34").Value = 2483428.5
$ws.Range("I34").Value = 1390429.6
$ws.Range("J34").Value = 5215926
$ws.Range("K34").Value = 1390429.6
$ws.Range("L34").Value = 5215926
$ws.Range("M34").Value = -1390227.6
$ws.Range("N34").Value = -5216330
# row 58
$ws.Range("H58").Value = 3869012.5
$ws.Range("I58").Value = 2236670
$ws.Range("J58").Value = 9092509
$ws.Range("K58").Value = 2236670
$ws.Range("L58").Value = 9092509
$ws.Range("M58").Value = -2236467
$ws.Range("N58").Value = -9092915
# row 122
$ws.Range("H122").Value = 3094.9211
$ws.Range("I122").Value = 4866.591
$ws.Range("J122").Value = 658.875
$ws.Range("K122").Value = 14599.773
$ws.Range("L122").Value = 1976.625
$ws.Range("M122").Value = -12149.773
$ws.Range("N122").Value = -6876.625
# row 132
$ws.Range("H132").Value = 2176953.8
$ws.Range("I132").Value = 4168952.2
$ws.Range("J132").Value = 3864.6365
$ws.Range("K132").Value = 12506856.6
$ws.Range("L132").Value = 11593.9095
$ws.Range("M132").Value = -12504326.6
$ws.Range("N132").Value = -16653.9095
# row 134
$ws.Range("H134").Value = 1825578.9
$ws.Range("I134").Value = 8653.714
$ws.Range("J134").Value = 5005198
$ws.Range("K134").Value = 25961.142
$ws.Range("L134").Value = 15015594
$ws.Range("M134").Value = -23426.142
$ws.Range("N134").Value = -15020664
# row 136
$ws.Range("H136").Value = 3869012.5
$ws.Range("I136").Value = 2236670
$ws.Range("J136").Value = 9092509
$ws.Range("K136").Value = 6710010
$ws.Range("L136").Value = 27277527
$ws.Range("M136").Value = -6707460
$ws.Range("N136").Value = -27282627

$ws = $wb.Worksheets.Item("CUL")
# row 3
$ws.Range("H3").Value = 6773.0625
$ws.Range("I3").Value = 5566.846
$ws.Range("J3").Value = 12000
$ws.Range("K3").Value = 16700.538
$ws.Range("L3").Value = 36000
$ws.Range("M3").Value = -16588.538
$ws.Range("N3").Value = -36224
# row 131
$ws.Range("H131").Value = 983.25
$ws.Range("I131").Value = 426.2
$ws.Range("J131").Value = 1168.9333
$ws.Range("K131").Value = 1278.6
$ws.Range("L131").Value = 3506.7999
$ws.Range("M131").Value = 3761.4
$ws.Range("N131").Value = -13586.7999

$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value = 13022096
$ws.Range("I132").Value = 16508847
$ws.Range("J132").Value = 8267435
$ws.Range("K132").Value = 49526541
$ws.Range("L132").Value = 24802305
$ws.Range("M132").Value = -49524011
$ws.Range("N132").Value = -24807365

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 5679.5713
$ws.Range("I22").Value = 4760
$ws.Range("J22").Value = 6190.4443
$ws.Range("K22").Value = 4760
$ws.Range("L22").Value = 6190.4443
$ws.Range("M22").Value = -4465
$ws.Range("N22").Value = -6780.4443
# row 27
$ws.Range("H27").Value = 5679.5713
$ws.Range("I27").Value = 4760
$ws.Range("J27").Value = 6190.4443
$ws.Range("K27").Value = 4760
$ws.Range("L27").Value = 6190.4443
$ws.Range("M27").Value = -4653
$ws.Range("N27").Value = -6404.4443
# row 46
$ws.Range("H46").Value = 378.27274
$ws.Range("J46").Value = 352.33334
$ws.Range("L46").Value = 352.33334
$ws.Range("N46").Value = -728.33334
# row 136
$ws.Range("H136").Value = 16178561
$ws.Range("I136").Value = 36976548
$ws.Range("J136").Value = 2349.5557
$ws.Range("K136").Value = 110929644
$ws.Range("L136").Value = 7048.6671
$ws.Range("M136").Value = -110927094
$ws.Range("N136").Value = -12148.6671

$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 35733416
$ws.Range("I62").Value = 83369800
$ws.Range("J62").Value = 6125
$ws.Range("K62").Value = 83369800
$ws.Range("L62").Value = 6125
$ws.Range("M62").Value = -83369176
$ws.Range("N62").Value = -7373
# row 65
$ws.Range("H65").Value = 35733416
$ws.Range("I65").Value = 83369800
$ws.Range("J65").Value = 6125
$ws.Range("K65").Value = 416849000
$ws.Range("L65").Value = 30625
$ws.Range("M65").Value = -416845880
$ws.Range("N65").Value = -36865
# row 86
$ws.Range("H86").Value = 20850
$ws.Range("J86").Value = 20850
$ws.Range("L86").Value = 20850
$ws.Range("N86").Value = -23096
# row 89
$ws.Range("H89").Value = 20850
$ws.Range("J89").Value = 20850
$ws.Range("L89").Value = 104250
$ws.Range("N89").Value = -115482
# row 132
$ws.Range("H132").Value = 760103.75
$ws.Range("I132").Value = 2640
$ws.Range("J132").Value = 3032495
$ws.Range("K132").Value = 7920
$ws.Range("L132").Value = 9097485
$ws.Range("M132").Value = -5390
$ws.Range("N132").Value = -9102545
# row 136
$ws.Range("H136").Value = 11979.333
$ws.Range("I136").Value = 9758.23
$ws.Range("J136").Value = 15588.625
$ws.Range("K136").Value = 29274.69
$ws.Range("L136").Value = 46765.875
$ws.Range("M136").Value = -26724.69
$ws.Range("N136").Value = -51865.875
